$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text update: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Use Replace so every cell referencing this text (Overview E/F columns and the
# Status column on the language sheets) is updated consistently.
$wsOverview.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
$wsZhCn.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
$wsDeDe.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null

# --- Latest Handback DateTime: give each language its own real timestamp ---
$wsZhCn.Cells.Replace("0001-01-01 00:00:00", "2016-08-19 23:07:31") | Out-Null
$wsDeDe.Cells.Replace("0001-01-01 00:00:00", "2016-08-19 23:07:37") | Out-Null

# --- Latest Target File (I) / Latest Handback File (J) population ---
$wsZhCn.Range("I2").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md"
$wsZhCn.Range("I3").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md"
$wsZhCn.Range("J2").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.352b1314ba9c16b46a9966aa39ddcfca17a91a00.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.352b1314ba9c16b46a9966aa39ddcfca17a91a00.zh-cn.xlf"

$wsDeDe.Range("I2").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md"
$wsDeDe.Range("I3").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md"
$wsDeDe.Range("J2").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.352b1314ba9c16b46a9966aa39ddcfca17a91a00.de-de.xlf"
$wsDeDe.Range("J3").Value = "b026ec3f-2496-45ab-bd83-7ca6f39dba91.352b1314ba9c16b46a9966aa39ddcfca17a91a00.de-de.xlf"

# --- Hyperlinks for the newly-populated "Latest Target File" cells ---
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f30cd4db1174732741ebe9b3093ffe1e1c85b451/e2e/b026ec3f-2496-45ab-bd83-7ca6f39dba91.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, "", "", "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetUrl, "", "", "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, "", "", "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetUrl, "", "", "b026ec3f-2496-45ab-bd83-7ca6f39dba91.md") | Out-Null

# --- Column width adjustments (values widened to fit the longer content) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

Write-Host "done"
